$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1) Highlight (gray fill) a set of existing rows' TestID cell (column A) --
#    mirrors the "open data provider run in parallel" review pass that
#    marked these rows with a light-gray background.
# ---------------------------------------------------------------------------
$ws.Range("A16").Interior.Color = 14277081
$ws.Range("A66:A69").Interior.Color = 14277081
$ws.Range("A76:A114").Interior.Color = 14277081

# ---------------------------------------------------------------------------
# 2) Append three new test cases (insert_021 / insert_022 / insert_023)
#    right after the last existing row (119), reusing row 119's exact
#    per-column formatting by copying it down before overwriting values.
# ---------------------------------------------------------------------------
$ws.Rows(119).Copy()
$ws.Rows(120).Insert(-4121)
$ws.Rows(119).Copy()
$ws.Rows(121).Insert(-4121)
$ws.Rows(119).Copy()
$ws.Rows(122).Insert(-4121)

# -- row 120: insert_021 ------------------------------------------------
$ws.Range("A120").Value = "insert_021"
$ws.Range("C120").Value = "指定replica为1,插入数据"
$ws.Range("D120").Value = "insert"
$ws.Range("F120").Value = "schema28"
$ws.Range("G120").Value = "insert_value17"
$ws.Range("H120").Value = "8"
$ws.Range("I120").Value = "select * from `$schema28"
$ws.Range("J120").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/insert/expectedresult/insert_021.csv"

# -- row 121: insert_022 ------------------------------------------------
$ws.Range("A121").Value = "insert_022"
$ws.Range("C121").Value = "指定replica为2,插入数据"
$ws.Range("D121").Value = "insert"
$ws.Range("F121").Value = "schema29"
$ws.Range("G121").Value = "insert_value18"
$ws.Range("H121").Value = "3"
$ws.Range("I121").Value = "select * from `$schema29"
$ws.Range("J121").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/insert/expectedresult/insert_022.csv"

# -- row 122: insert_023 ------------------------------------------------
$ws.Range("A122").Value = "insert_023"
$ws.Range("C122").Value = "指定replica为3,插入数据"
$ws.Range("D122").Value = "insert"
$ws.Range("F122").Value = "schema30"
$ws.Range("G122").Value = "insert_value19"
$ws.Range("H122").Value = "3"
$ws.Range("I122").Value = "select * from `$schema30"
$ws.Range("J122").Value = "src/test/resources/io.dingodb.test/testdata/mysqlcases/dml/insert/expectedresult/insert_023.csv"

# ---------------------------------------------------------------------------
# 3) Restore the view state (scroll position / active selection) to match
#    where the author ended up after adding the new rows.
# ---------------------------------------------------------------------------
$ws.Range("A85").Select()
$excel.ActiveWindow.ScrollRow = 85
$ws.Range("I116").Select()
